$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = -16.73917251691514
$ws.Range("F2").Value = 3.103248084001659
$ws.Range("G2").Value = -10.33120499552526
$ws.Range("E3").Value = -17.27966689307319
$ws.Range("F3").Value = 3.016918000745675
$ws.Range("G3").Value = -10.23750717892083
$ws.Range("E4").Value = -17.80771874675494
$ws.Range("F4").Value = 3.095415895000861
$ws.Range("G4").Value = -10.2840309682664
$ws.Range("E5").Value = -18.21477256948742
$ws.Range("F5").Value = 3.465963503458079
$ws.Range("G5").Value = -9.960344490628557
$ws.Range("E6").Value = -18.89835350961946
$ws.Range("F6").Value = 3.567444618994632
$ws.Range("G6").Value = -9.657685631525322
$ws.Range("E7").Value = -19.51756578534171
$ws.Range("F7").Value = 3.566603709813523
$ws.Range("G7").Value = -9.378161552916326
$ws.Range("E8").Value = -20.06695326499068
$ws.Range("F8").Value = 3.833812380124759
$ws.Range("G8").Value = -9.184708440199369
$ws.Range("E9").Value = -20.54137271334237
$ws.Range("F9").Value = 4.154492118539444
$ws.Range("G9").Value = -9.363973654988538
$ws.Range("E10").Value = -21.18314286674856
$ws.Range("F10").Value = 4.475895429970433
$ws.Range("G10").Value = -8.906460392382638
$ws.Range("E11").Value = -22.16304162407745
$ws.Range("F11").Value = 4.599215739181496
$ws.Range("G11").Value = -8.698946495916433
$ws.Range("E12").Value = -22.97715415854926
$ws.Range("F12").Value = 4.967739298795811
$ws.Range("G12").Value = -8.245530221065007
$ws.Range("E13").Value = -23.6132726200173
$ws.Range("F13").Value = 5.780072234768064
$ws.Range("G13").Value = -7.944171837788369
$ws.Range("E14").Value = -24.79367443796525
$ws.Range("F14").Value = 6.097554562691787
$ws.Range("G14").Value = -7.105027588157954
$ws.Range("E15").Value = -25.33368969145034
$ws.Range("F15").Value = 6.538909657602534
$ws.Range("G15").Value = -6.618527403838115
$ws.Range("E16").Value = -26.24358275945188
$ws.Range("F16").Value = 6.757286927327024
$ws.Range("G16").Value = -6.445495672804261
$ws.Range("E17").Value = -27.3119554295479
$ws.Range("F17").Value = 7.200050056215442
$ws.Range("G17").Value = -6.112290298796539
$ws.Range("E18").Value = -27.92894076264227
$ws.Range("F18").Value = 7.524293880636109
$ws.Range("G18").Value = -6.049285667302607
$ws.Range("E19").Value = -28.98439423209247
$ws.Range("F19").Value = 7.79429906287522
$ws.Range("G19").Value = -5.738207938374543
$ws.Range("E20").Value = -29.67376131185151
$ws.Range("F20").Value = 8.115482368997197
$ws.Range("G20").Value = -5.635436125025123
$ws.Range("E21").Value = -30.36263704642043
$ws.Range("F21").Value = 8.402119952597555
$ws.Range("G21").Value = -5.59848501112498
$ws.Range("E22").Value = -30.66216205212197
$ws.Range("F22").Value = 8.417994557894428
$ws.Range("G22").Value = -5.653481449370906
$ws.Range("E23").Value = -31.01262073136304
$ws.Range("F23").Value = 8.532255537381097
$ws.Range("G23").Value = -5.464413775813571
$ws.Range("E24").Value = -31.13014267842994
$ws.Range("F24").Value = 8.963392607939991
$ws.Range("G24").Value = -5.509920651730816
$ws.Range("E25").Value = -31.18684538007242
$ws.Range("F25").Value = 8.87290686884713
$ws.Range("G25").Value = -5.426567973656783
$ws.Range("E26").Value = -31.10761169028376
$ws.Range("F26").Value = 8.476183517624449
$ws.Range("G26").Value = -5.419752698084303
$ws.Range("E27").Value = -31.2751995121713
$ws.Range("F27").Value = 8.546966559043645
$ws.Range("G27").Value = -5.079790716586506
$ws.Range("E28").Value = -30.96217351400677
$ws.Range("F28").Value = 8.45078023794396
$ws.Range("G28").Value = -5.312976788110881
$ws.Range("E29").Value = -30.86274333685089
$ws.Range("F29").Value = 8.063267775658547
$ws.Range("G29").Value = -5.000227018834332
$ws.Range("E30").Value = -31.03710243324911
$ws.Range("F30").Value = 7.847760352964936
$ws.Range("G30").Value = -5.121899732731244
$ws.Range("E31").Value = -30.66291984818634
$ws.Range("F31").Value = 7.734203390467102
$ws.Range("G31").Value = -5.092179459987269
$ws.Range("E32").Value = -30.72304485463566
$ws.Range("F32").Value = 7.500288856919556
$ws.Range("G32").Value = -5.213373051211223
$ws.Range("E33").Value = -30.2984346082441
$ws.Range("F33").Value = 7.458522071255502
$ws.Range("G33").Value = -4.906103858632485
$ws.Range("E34").Value = -30.03612472281351
$ws.Range("F34").Value = 7.390530652757316
$ws.Range("G34").Value = -4.870409219497139
$ws.Range("E35").Value = -29.62258563247208
$ws.Range("F35").Value = 7.447844480258159
$ws.Range("G35").Value = -4.965974636724725
$ws.Range("E36").Value = -29.34637385601542
$ws.Range("F36").Value = 7.425193711443742
$ws.Range("G36").Value = -4.918741941383461
$ws.Range("E37").Value = -29.15831820687954
$ws.Range("F37").Value = 7.428102670529557
$ws.Range("G37").Value = -4.798403926361216
$ws.Range("E38").Value = -28.55720747908204
$ws.Range("F38").Value = 7.138575683870854
$ws.Range("G38").Value = -4.696776140618654
$ws.Range("E39").Value = -28.04287662317229
$ws.Range("F39").Value = 7.128739002054621
$ws.Range("G39").Value = -4.677298337260866
$ws.Range("E40").Value = -27.64174583225908
$ws.Range("F40").Value = 7.11946944503495
$ws.Range("G40").Value = -4.473157855532832
$ws.Range("E41").Value = -27.47574693760329
$ws.Range("F41").Value = 6.928265275476774
$ws.Range("G41").Value = -4.688831504459918
$ws.Range("E42").Value = -26.90051127864894
$ws.Range("F42").Value = 6.785804504381739
$ws.Range("G42").Value = -4.482334521422032
$ws.Range("E43").Value = -26.500013614125
$ws.Range("F43").Value = 6.797401228670061
$ws.Range("G43").Value = -4.460363324562116
$ws.Range("E44").Value = -25.84442468280947
$ws.Range("F44").Value = 6.831238045195979
$ws.Range("G44").Value = -4.504149270062205
$ws.Range("E45").Value = -25.38822656305078
$ws.Range("F45").Value = 6.781477733304519
$ws.Range("G45").Value = -4.463971411629899
$ws.Range("E46").Value = -24.83691770370311
$ws.Range("F46").Value = 6.703864749292245
$ws.Range("G46").Value = -4.430569716715135
$ws.Range("E47").Value = -24.76217456672171
$ws.Range("F47").Value = 6.880093890817061
$ws.Range("G47").Value = -4.519794092036333
$ws.Range("E48").Value = -24.35641633080898
$ws.Range("F48").Value = 6.80507208044425
$ws.Range("G48").Value = -4.637531155405379
$ws.Range("E49").Value = -23.64784767658012
$ws.Range("F49").Value = 6.980597204980231
$ws.Range("G49").Value = -4.758714968615601
$ws.Range("E50").Value = -23.43367495276118
$ws.Range("F50").Value = 6.868076711938183
$ws.Range("G50").Value = -4.956685523677587
$ws.Range("E51").Value = -22.99567860556801
$ws.Range("F51").Value = 6.879135645471146
$ws.Range("G51").Value = -4.70895465672414
$ws.Range("E52").Value = -22.36690343241407
$ws.Range("F52").Value = 6.812459369820159
$ws.Range("G52").Value = -4.888884776447211
$ws.Range("E53").Value = -21.95146985191171
$ws.Range("F53").Value = 7.15820015743465
$ws.Range("G53").Value = -4.918228595662434
$ws.Range("E54").Value = -21.33269758680749
$ws.Range("F54").Value = 7.173561417010498
$ws.Range("G54").Value = -4.808885957083882
$ws.Range("E55").Value = -21.09269623971063
$ws.Range("F55").Value = 7.091802555175079
$ws.Range("G55").Value = -5.180259807701608
$ws.Range("E56").Value = -20.56007805361518
$ws.Range("F56").Value = 6.990605002036806
$ws.Range("G56").Value = -5.354548013500259
$ws.Range("E57").Value = -20.14318754906648
$ws.Range("F57").Value = 7.036194991070788
$ws.Range("G57").Value = -5.508312168471601
$ws.Range("E58").Value = -19.86608841786347
$ws.Range("F58").Value = 6.774378895333759
$ws.Range("G58").Value = -5.904404837808449
$ws.Range("E59").Value = -19.35291870108462
$ws.Range("F59").Value = 6.952030737856847
$ws.Range("G59").Value = -5.875422805101377
$ws.Range("E60").Value = -18.81773378638404
$ws.Range("F60").Value = 6.777190074282235
$ws.Range("G60").Value = -6.246219752908808
$ws.Range("E61").Value = -18.59032163196948
$ws.Range("F61").Value = 6.582446263752412
$ws.Range("G61").Value = -6.271437250328355
$ws.Range("E62").Value = -18.35154742559623
$ws.Range("F62").Value = 6.549997925176697
$ws.Range("G62").Value = -6.353235224218709
$ws.Range("E63").Value = -17.69447468069659
$ws.Range("F63").Value = 6.554999379201552
$ws.Range("G63").Value = -6.277954296481953
$ws.Range("E64").Value = -17.74024556298465
$ws.Range("F64").Value = 6.336030539646164
$ws.Range("G64").Value = -6.617842942876747
$ws.Range("E65").Value = -17.37410294971452
$ws.Range("F65").Value = 6.440938848996425
$ws.Range("G65").Value = -6.472116315194613
$ws.Range("E66").Value = -17.20844384103598
$ws.Range("F66").Value = 6.458074818064962
$ws.Range("G66").Value = -6.553792063913294
$ws.Range("E67").Value = -17.05597527288437
$ws.Range("F67").Value = 6.09878659242225
$ws.Range("G67").Value = -6.567021716495166
$ws.Range("E68").Value = -16.7324696885006
$ws.Range("F68").Value = 6.080081252149433
$ws.Range("G68").Value = -6.69778798316454
$ws.Range("E69").Value = -16.76170106055789
$ws.Range("F69").Value = 6.165413978004566
$ws.Range("G69").Value = -6.410226377266337
$ws.Range("E70").Value = -16.41227885077261
$ws.Range("F70").Value = 6.115252767550019
$ws.Range("G70").Value = -6.256912010926751
$ws.Range("E71").Value = -16.5912605031635
$ws.Range("F71").Value = 5.908638448347326
$ws.Range("G71").Value = -5.997628420746784
$ws.Range("E72").Value = -16.56600878269589
$ws.Range("F72").Value = 5.929196722222704
$ws.Range("G72").Value = -5.973638064050832
$ws.Range("E73").Value = -16.40793741267479
$ws.Range("F73").Value = 5.904482792510449
$ws.Range("G73").Value = -6.077876579460326
$ws.Range("E74").Value = -16.5453038386145
$ws.Range("F74").Value = 5.777173053695984
$ws.Range("G74").Value = -5.654899261362312
$ws.Range("E75").Value = -16.71340256171964
$ws.Range("F75").Value = 5.72130148322088
$ws.Range("G75").Value = -5.647795534384684
$ws.Range("E76").Value = -16.81079646751545
$ws.Range("F76").Value = 5.592305037037328
$ws.Range("G76").Value = -5.480733280735334
$ws.Range("E77").Value = -16.98164281247979
$ws.Range("F77").Value = 5.59686648044416
$ws.Range("G77").Value = -5.24358222464189
$ws.Range("E78").Value = -17.32831740040587
$ws.Range("F78").Value = 5.364213310668279
$ws.Range("G78").Value = -4.879390325111662
$ws.Range("E79").Value = -17.56636803376282
$ws.Range("F79").Value = 5.498450772213163
$ws.Range("G79").Value = -4.543310215066194
$ws.Range("E80").Value = -17.94757856619678
$ws.Range("F80").Value = 5.318989997149315
$ws.Range("G80").Value = -4.240842027230768
$ws.Range("E81").Value = -18.69960560246563
$ws.Range("F81").Value = 5.294955639391562
$ws.Range("G81").Value = -3.989146175708252
$ws.Range("E82").Value = -19.386848408741
$ws.Range("F82").Value = 5.116091323165478
$ws.Range("G82").Value = -3.726186052364365
$ws.Range("E83").Value = -20.028662563209
$ws.Range("F83").Value = 4.863847902873868
$ws.Range("G83").Value = -3.584659081580911
$ws.Range("E84").Value = -21.18079614345244
$ws.Range("F84").Value = 4.606563916502469
$ws.Range("G84").Value = -3.135921586301122
$ws.Range("E85").Value = -22.01180997466318
$ws.Range("F85").Value = 4.805355824718098
$ws.Range("G85").Value = -3.103844812247293
$ws.Range("E86").Value = -22.8427358037503
$ws.Range("F86").Value = 4.468097464017704
$ws.Range("G86").Value = -2.807610108167183
$ws.Range("E87").Value = -23.5850630503952
$ws.Range("F87").Value = 4.31796584114848
$ws.Range("G87").Value = -2.731449159194382
$ws.Range("E88").Value = -24.85820932860852
$ws.Range("F88").Value = 4.281185842488679
$ws.Range("G88").Value = -2.637076659642323
$ws.Range("E89").Value = -26.10451007011562
$ws.Range("F89").Value = 3.71836337196943
$ws.Range("G89").Value = -2.554765340030942
$ws.Range("E90").Value = -27.50126999795201
$ws.Range("F90").Value = 3.811650511997034
$ws.Range("G90").Value = -2.55411999112451
$ws.Range("E91").Value = -29.19023364370675
$ws.Range("F91").Value = 3.427697246710735
$ws.Range("G91").Value = -2.746565968426883
$ws.Range("E92").Value = -30.88354850557305
$ws.Range("F92").Value = 3.188971930406145
$ws.Range("G92").Value = -2.918287445620404
$ws.Range("E93").Value = -32.43171608707094
$ws.Range("F93").Value = 2.847724140102349
$ws.Range("G93").Value = -3.271703974015951
$ws.Range("E94").Value = -34.40046537330205
$ws.Range("F94").Value = 2.340670570913999
$ws.Range("G94").Value = -3.143924890542262
$ws.Range("E95").Value = -36.60736551753089
$ws.Range("F95").Value = 1.842260765866355
$ws.Range("G95").Value = -3.276597869889733
$ws.Range("E96").Value = -38.41338665470018
$ws.Range("F96").Value = 1.576020118914767
$ws.Range("G96").Value = -3.590990345473565
$ws.Range("E97").Value = -40.52297111724307
$ws.Range("F97").Value = 1.108811955691774
$ws.Range("G97").Value = -3.840882153462189
$ws.Range("E98").Value = -42.29232959191736
$ws.Range("F98").Value = 0.7353113870869397
$ws.Range("G98").Value = -4.16587888494036
$ws.Range("E99").Value = -44.80218438016489
$ws.Range("F99").Value = 0.3961072016535136
$ws.Range("G99").Value = -4.346009453944975
$ws.Range("E100").Value = -47.32341588739173
$ws.Range("F100").Value = -0.3613515212444862
$ws.Range("G100").Value = -4.679742840694323
$ws.Range("E101").Value = -49.19887803359524
$ws.Range("F101").Value = -0.3484445431158308
$ws.Range("G101").Value = -5.092653693653359
$ws.Range("E102").Value = -51.38210193981932
$ws.Range("F102").Value = -0.7973140415810259
$ws.Range("G102").Value = -5.426587529684251
